$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the differing values between row 2 and row 3 (use .Value2 for reads,
# since it reliably returns the underlying scalar for both numbers and text)

# Column D (Fecha)
$d2 = $ws.Range("D2").Value2
$d3 = $ws.Range("D3").Value2
$ws.Range("D2").Value = $d3
$ws.Range("D3").Value = $d2

# Column M (Volumen)
$m2 = $ws.Range("M2").Value2
$m3 = $ws.Range("M3").Value2
$ws.Range("M2").Value = $m3
$ws.Range("M3").Value = $m2

# Column N (Precio mínimo)
$n2 = $ws.Range("N2").Value2
$n3 = $ws.Range("N3").Value2
$ws.Range("N2").Value = $n3
$ws.Range("N3").Value = $n2

# Column O (Precio máximo)
$o2 = $ws.Range("O2").Value2
$o3 = $ws.Range("O3").Value2
$ws.Range("O2").Value = $o3
$ws.Range("O3").Value = $o2

# Column P (Precio promedio ponderado)
$p2 = $ws.Range("P2").Value2
$p3 = $ws.Range("P3").Value2
$ws.Range("P2").Value = $p3
$ws.Range("P3").Value = $p2

# Column R (Origen)
$r2 = $ws.Range("R2").Value2
$r3 = $ws.Range("R3").Value2
$ws.Range("R2").Value = $r3
$ws.Range("R3").Value = $r2

# Column S (Precio $/Kg)
$s2 = $ws.Range("S2").Value2
$s3 = $ws.Range("S3").Value2
$ws.Range("S2").Value = $s3
$ws.Range("S3").Value = $s2
